$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Athena README tasks: flip the three "Open" status cells to "Done"/"done"
# (row 4 used a lowercase "done" string).
$ws.Range("E3").Value = "Done"
$ws.Range("E4").Value = "done"
$ws.Range("E22").Value = "Done"

# Widen the Section column (A) so the longer section names fit; this also
# clears Excel's bestFit auto-size flag since the width is now explicit.
$ws.Columns("A").ColumnWidth = 24.8

# Scroll/selection moved to A24 before saving.
$ws.Range("A24").Select() | Out-Null
